# Updates for 4.0.1 - Cap Ret per Unit Net Loss
# - Reworks retirement mechanism/data to work in terms of MW/($/MW) instead of
#   Fraction of start year capacity, i.e. MW/($/MWh)
# - Updates economic retirement parameters for coal and BAU retirements (111 rules)
# - Updates "About" sheet narrative text (natural gas peakers removed from the
#   "low retirement fraction" group; new note about biomass/CHP colocation)

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("CRpUNL")

# ---------------------------------------------------------------------------
# CRpUNL sheet: rename headers and switch data from "Fraction of start year
# capacity" to "MW retired", re-expressing the unit accordingly.
# ---------------------------------------------------------------------------
$wsData.Range("B1").Value = "MW retired"
$wsData.Range("A1").Value = "Unit: MW/(`$/MW)"

# ---------------------------------------------------------------------------
# About sheet: revise the "low retirement" explanatory sentence and append a
# new explanatory note about biomass / CHP colocation.
# ---------------------------------------------------------------------------
$wsAbout.Range("A10").Value = "These includes: natural gas steam turbines and petroleum plants. For these plant types we set the "
$wsAbout.Range("A16").Value = "Likewise, biomass plants are often colocated with cheap supply and part of integrated"
$wsAbout.Range("A17").Value = "CHP or industrial systems, and we therefore do not subject them to economic retirement."

# Plant types that keep a normal (non-reliability-protected) retirement rate
$wsData.Range("B2").Value  = 0.03   # hard coal
$wsData.Range("B3").Value  = 0.03   # natural gas steam turbine
$wsData.Range("B4").Value  = 0.03   # natural gas combined cycle
$wsData.Range("B5").Value  = 0.03   # nuclear
# B6 hydro stays at 0 (gray-filled, unchanged)
$wsData.Range("B7").Value  = 0.03   # onshore wind
$wsData.Range("B8").Value  = 0.03   # solar PV

# Plant types newly set to zero and flagged with the reliability-exempt fill
# (solar thermal, biomass, geothermal, petroleum)
$wsData.Range("B6").Copy()
$wsData.Range("B9:B12").PasteSpecial(-4122)   # xlPasteFormats (copy the existing gray-fill style)
$excel.CutCopyMode = $false

$wsData.Range("B9").Value  = 0   # solar thermal
$wsData.Range("B10").Value = 0   # biomass
$wsData.Range("B11").Value = 0   # geothermal
$wsData.Range("B12").Value = 0   # petroleum

$wsData.Range("B13").Value = 0.03   # natural gas peaker
$wsData.Range("B14").Value = 0.03   # lignite
$wsData.Range("B15").Value = 0.03   # offshore wind
# B16 crude oil stays at 0 (gray-filled, unchanged)
# B17 heavy or residual fuel oil stays at 0 (gray-filled, unchanged)

$wsData.Range("B6").Copy()
$wsData.Range("B18").PasteSpecial(-4122)   # xlPasteFormats (copy the existing gray-fill style)
$excel.CutCopyMode = $false
$wsData.Range("B18").Value = 0   # municipal solid waste

$wsData.Range("B19").Value = 0.03   # hard coal w CCS
$wsData.Range("B20").Value = 0.03   # natural gas combined cycle w CCS
$wsData.Range("B21").Value = 0.03   # biomass w CCS
$wsData.Range("B22").Value = 0.03   # lignite w CCS
$wsData.Range("B23").Value = 0.03   # small modular reactor
$wsData.Range("B24").Value = 0.03   # hydrogen combustion turbine
$wsData.Range("B25").Value = 0.03   # hydrogen combined cycle

# ---------------------------------------------------------------------------
# Restore selections to match the final authored state.
# ---------------------------------------------------------------------------
$wsData.Range("D14").Select() | Out-Null
$wsAbout.Activate() | Out-Null
$wsAbout.Range("A18").Select() | Out-Null
